# "Update points 09876543 -> 0.00"
#
# Row 16 held phone "09876543" (stored as text, leading zero preserved) with
# 0 points. A new entry for phone 9876543 (entered as a number, so the
# leading zero is dropped) with 0 points is added directly above it - the
# original text row is pushed down to row 17 and is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 16; this shifts the old row 16 (phone "09876543",
# blank birthday, 0 points) down to row 17 unchanged.
$ws.Rows.Item(16).Insert()

# Populate the new row 16: phone entered as a plain number (9876543) and
# 0 points; birthday stays blank, same as the row it was inserted above.
$ws.Range("A16").Value = 9876543
$ws.Range("C16").Value = 0
